$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting existing rows 72-159 down to 73-160.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly data point.
$ws.Cells.Item(72, 1).Value = 4
$ws.Cells.Item(72, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(72, 3).Value = "Los Lagos"
$ws.Cells.Item(72, 4).Value = 45128
$ws.Cells.Item(72, 5).Value = 10
$ws.Cells.Item(72, 6).Value = 100112031
$ws.Cells.Item(72, 7).Value = "Poroto verde"
$ws.Cells.Item(72, 8).Value = "Magnum"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 45
$ws.Cells.Item(72, 11).Value = 28000
$ws.Cells.Item(72, 12).Value = 28000
$ws.Cells.Item(72, 13).Value = 28000
$ws.Cells.Item(72, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(72, 15).Value = "Perú"
$ws.Cells.Item(72, 16).Value = 1120
$ws.Cells.Item(72, 17).Value = 25
$ws.Cells.Item(72, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D.
$ws.Cells.Item(72, 4).NumberFormat = $ws.Cells.Item(73, 4).NumberFormat
